$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3409
$ws.Range("J40").Value = 3409
$ws.Range("L40").Value = 3409
$ws.Range("N40").Value = -3759
$ws.Range("H62").Value = 2345.842
$ws.Range("I62").Value = 2387.3333
$ws.Range("K62").Value = 2387.3333
$ws.Range("M62").Value = -1763.3333
$ws.Range("H65").Value = 2345.842
$ws.Range("I65").Value = 2387.3333
$ws.Range("K65").Value = 11936.6665
$ws.Range("M65").Value = -8816.666499999999
$ws.Range("H69").Value = 9230.462
$ws.Range("J69").Value = 11833.333
$ws.Range("L69").Value = 35499.999
$ws.Range("N69").Value = -37247.999
$ws.Range("H70").Value = 51508.05
$ws.Range("I70").Value = 101037
$ws.Range("J70").Value = 1979.1
$ws.Range("K70").Value = 303111
$ws.Range("L70").Value = 5937.299999999999
$ws.Range("M70").Value = -302841
$ws.Range("N70").Value = -6477.299999999999
$ws.Range("H72").Value = 9230.462
$ws.Range("J72").Value = 11833.333
$ws.Range("L72").Value = 106499.997
$ws.Range("N72").Value = -115235.997
$ws.Range("H73").Value = 51508.05
$ws.Range("I73").Value = 101037
$ws.Range("J73").Value = 1979.1
$ws.Range("K73").Value = 303111
$ws.Range("L73").Value = 5937.299999999999
$ws.Range("M73").Value = -302175
$ws.Range("N73").Value = -7809.299999999999
$ws.Range("H76").Value = 4499.5
$ws.Range("I76").Value = 4666
$ws.Range("K76").Value = 4666
$ws.Range("M76").Value = -4351
$ws.Range("H79").Value = 4499.5
$ws.Range("I79").Value = 4666
$ws.Range("K79").Value = 4666
$ws.Range("M79").Value = -3574
$ws.Range("H82").Value = 1818.3077
$ws.Range("I82").Value = 1849.4166
$ws.Range("K82").Value = 5548.2498
$ws.Range("M82").Value = -5142.2498
$ws.Range("H85").Value = 1818.3077
$ws.Range("I85").Value = 1849.4166
$ws.Range("K85").Value = 5548.2498
$ws.Range("M85").Value = -4144.2498
$ws.Range("H92").Value = 1803626.1
$ws.Range("I92").Value = 976811
$ws.Range("K92").Value = 976811
$ws.Range("M92").Value = -975563
$ws.Range("H111").Value = 2859.7778
$ws.Range("I111").Value = 3277.8
$ws.Range("K111").Value = 9833.400000000001
$ws.Range("M111").Value = -6766.400000000001
$ws.Range("H113").Value = 2882.1667
$ws.Range("I113").Value = 2858.6
$ws.Range("K113").Value = 2858.6
$ws.Range("M113").Value = 395.4000000000001
$ws.Range("H116").Value = 3479.182
$ws.Range("I116").Value = 3479.182
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 3479.182
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -37.18199999999979
$ws.Range("H130").Value = 121497.5
$ws.Range("J130").Value = 121497.5
$ws.Range("L130").Value = 121497.5
$ws.Range("N130").Value = -131537.5
$ws.Range("H137").Value = 3450055.5
$ws.Range("I137").Value = 1103.7778
$ws.Range("J137").Value = 5002083.5
$ws.Range("K137").Value = 3311.3334
$ws.Range("L137").Value = 15006250.5
$ws.Range("M137").Value = -761.3334000000004
$ws.Range("N137").Value = -15011350.5
$ws.Range("N116").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 486632.2
$ws.Range("I32").Value = 505963.1
$ws.Range("J32").Value = 100014
$ws.Range("K32").Value = 505963.1
$ws.Range("L32").Value = 100014
$ws.Range("M32").Value = -505676.1
$ws.Range("N32").Value = -100588

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 25717814
$ws.Range("I134").Value = 2716.0435
$ws.Range("K134").Value = 8148.130500000001
$ws.Range("M134").Value = -5613.130500000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2852.077
$ws.Range("I31").Value = 1611.9286
$ws.Range("K31").Value = 1611.9286
$ws.Range("M31").Value = -1316.9286
$ws.Range("H32").Value = 1003500
$ws.Range("I32").Value = 1003500
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1003500
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1003184
$ws.Range("H34").Value = 2852.077
$ws.Range("I34").Value = 1611.9286
$ws.Range("K34").Value = 1611.9286
$ws.Range("M34").Value = -1409.9286
$ws.Range("H62").Value = 336933
$ws.Range("I62").Value = 5400
$ws.Range("K62").Value = 5400
$ws.Range("M62").Value = -4776
$ws.Range("H65").Value = 336933
$ws.Range("I65").Value = 5400
$ws.Range("K65").Value = 27000
$ws.Range("M65").Value = -23880
$ws.Range("H107").Value = 1633.44
$ws.Range("I107").Value = 1335.3334
$ws.Range("J107").Value = 3198.5
$ws.Range("K107").Value = 1335.3334
$ws.Range("L107").Value = 3198.5
$ws.Range("M107").Value = 584.6666
$ws.Range("N107").Value = -7038.5
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("H134").Value = 2550
$ws.Range("I134").Value = 2019.7
$ws.Range("J134").Value = 3433.8333
$ws.Range("K134").Value = 6059.1
$ws.Range("L134").Value = 10301.4999
$ws.Range("M134").Value = -3524.1
$ws.Range("N134").Value = -15371.4999
$ws.Range("N32").ClearContents()
$ws.Range("N123").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 446.1875
$ws.Range("I5").Value = 451.35715
$ws.Range("J5").Value = 410
$ws.Range("K5").Value = 1354.07145
$ws.Range("L5").Value = 1230
$ws.Range("M5").Value = -1242.07145
$ws.Range("N5").Value = -1454
$ws.Range("H11").Value = 2111721
$ws.Range("I11").Value = 634.7273
$ws.Range("K11").Value = 1904.1819
$ws.Range("M11").Value = -1764.1819
$ws.Range("H12").Value = 28644.611
$ws.Range("J12").Value = 85835.836
$ws.Range("L12").Value = 257507.508
$ws.Range("N12").Value = -257853.508
$ws.Range("H48").Value = 5600
$ws.Range("J48").Value = 6200
$ws.Range("L48").Value = 18600
$ws.Range("N48").Value = -19100
$ws.Range("H131").Value = 5078330.5
$ws.Range("J131").Value = 3619262.8
$ws.Range("L131").Value = 10857788.4
$ws.Range("N131").Value = -10867868.4
$ws.Range("H134").Value = 27780876
$ws.Range("I134").Value = 31252862
$ws.Range("K134").Value = 93758586
$ws.Range("M134").Value = -93753516
$ws.Range("H135").Value = 446.1875
$ws.Range("I135").Value = 451.35715
$ws.Range("J135").Value = 410
$ws.Range("K135").Value = 4062.21435
$ws.Range("L135").Value = 3690
$ws.Range("M135").Value = -1527.21435
$ws.Range("N135").Value = -8760

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 2478.0908
$ws.Range("I24").Value = 4950
$ws.Range("J24").Value = 1928.7778
$ws.Range("K24").Value = 4950
$ws.Range("L24").Value = 1928.7778
$ws.Range("M24").Value = -4777
$ws.Range("N24").Value = -2274.7778
$ws.Range("H70").Value = 9869.615
$ws.Range("I70").Value = 9339.241
$ws.Range("K70").Value = 9339.241
$ws.Range("M70").Value = -9069.241
$ws.Range("H73").Value = 9869.615
$ws.Range("I73").Value = 9339.241
$ws.Range("K73").Value = 9339.241
$ws.Range("M73").Value = -8403.241
$ws.Range("H132").Value = 8616601
$ws.Range("I132").Value = 2772.6956
$ws.Range("K132").Value = 8318.086800000001
$ws.Range("M132").Value = -5788.086800000001
$ws.Range("H141").Value = 62862.375
$ws.Range("J141").Value = 62862.375
$ws.Range("L141").Value = 62862.375
$ws.Range("N141").Value = -73222.375

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("H93").Value = 1705.4615
$ws.Range("I93").Value = 1659.6
$ws.Range("K93").Value = 1659.6
$ws.Range("M93").Value = -411.5999999999999
$ws.Range("H100").Value = 2999
$ws.Range("I100").Value = 3000
$ws.Range("J100").Value = 2998.6667
$ws.Range("K100").Value = 3000
$ws.Range("L100").Value = 2998.6667
$ws.Range("M100").Value = -2459
$ws.Range("N100").Value = -4080.6667
$ws.Range("H132").Value = 7551.2856
$ws.Range("I132").Value = 1485
$ws.Range("K132").Value = 4455
$ws.Range("M132").Value = -1925
$ws.Range("H136").Value = 2811.0833
$ws.Range("I136").Value = 1656
$ws.Range("J136").Value = 5121.25
$ws.Range("K136").Value = 4968
$ws.Range("L136").Value = 15363.75
$ws.Range("M136").Value = -2418
$ws.Range("N136").Value = -20463.75
$ws.Range("H140").Value = 94628.875
$ws.Range("J140").Value = 94628.875
$ws.Range("L140").Value = 94628.875
$ws.Range("N140").Value = -104988.875
$ws.Range("N69").ClearContents()
$ws.Range("N72").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 10000000
$ws.Range("I5").Value = 10000000
$ws.Range("K5").Value = 10000000
$ws.Range("M5").Value = -9999888
$ws.Range("H46").Value = 58355.6
$ws.Range("J46").Value = 58355.6
$ws.Range("L46").Value = 58355.6
$ws.Range("N46").Value = -58817.6
$ws.Range("H113").Value = 312.6
$ws.Range("J113").Value = 266.375
$ws.Range("L113").Value = 799.125
$ws.Range("N113").Value = -5139.125
$ws.Range("H125").Value = 105102.8
$ws.Range("J125").Value = 105102.8
$ws.Range("L125").Value = 105102.8
$ws.Range("N125").Value = -114942.8
$ws.Range("H134").Value = 58355.6
$ws.Range("J134").Value = 58355.6
$ws.Range("L134").Value = 175066.8
$ws.Range("N134").Value = -180136.8

